# Auto-generated script applying cryptos.xlsx price/volume update (2023-07-20)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.289.55"
$ws.Range("E2").Value = "'  +0.97%  "

$ws.Range("D3").Value = "'1.921.43"
$ws.Range("E3").Value = "'  +0.59%  "

$ws.Range("D5").Value = "'0.8164"
$ws.Range("E5").Value = "'  +3.11%  "

$ws.Range("D6").Value = "'244.21"

$ws.Range("D7").Value = "'1.000"
$ws.Range("E7").Value = "'  -0.01%  "

$ws.Range("D8").Value = "'0.3264"
$ws.Range("E8").Value = "'  +3.49%  "

$ws.Range("D9").Value = "'27.32"
$ws.Range("E9").Value = "'  +4.07%  "

$ws.Range("D10").Value = "'0.07286"
$ws.Range("E10").Value = "'  +5.51%  "

$ws.Range("D11").Value = "'0.7967"
$ws.Range("E11").Value = "'  +7.41%  "

$ws.Range("D12").Value = "'0.08113"
$ws.Range("E12").Value = "'  +1.40%  "

$ws.Range("D13").Value = "'1.941.67"
$ws.Range("E13").Value = "'  +1.67%  "

$ws.Range("D14").Value = "'5.418"
$ws.Range("E14").Value = "'  +4.44%  "

$ws.Range("D15").Value = "'94.33"
$ws.Range("E15").Value = "'  +1.44%  "

$ws.Range("D16").Value = "'30.291.72"
$ws.Range("E16").Value = "'  +0.99%  "

$ws.Range("D17").Value = "'14.27"

$ws.Range("D18").Value = "'6.078"
$ws.Range("E18").Value = "'  +3.69%  "

$ws.Range("D19").Value = "'250.59"

$ws.Range("D20").Value = "'0.000007874"
$ws.Range("E20").Value = "'  +1.76%  "

$ws.Range("D21").Value = "'2.181.65"
$ws.Range("E21").Value = "'  +1.02%  "

$ws.Range("E22").Value = "'  +0.04%  "

$ws.Range("D23").Value = "'8.039"
$ws.Range("E23").Value = "'  +17.68%  "

$ws.Range("E24").Value = "'  +0.07%  "

$ws.Range("D25").Value = "'0.1682"
$ws.Range("E25").Value = "'  +21.21%  "

$ws.Range("D26").Value = "'9.512"
$ws.Range("E26").Value = "'  +3.04%  "

$ws.Range("D27").Value = "'167.90"
$ws.Range("E27").Value = "'  -0.08%  "

$ws.Range("D28").Value = "'19.07"
$ws.Range("E28").Value = "'  +0.83%  "

$ws.Range("D29").Value = "'2.155"
$ws.Range("E29").Value = "'  +6.03%  "

$ws.Range("E30").Value = "'  +0.47%  "

$ws.Range("E31").Value = "'  +2.52%  "

$ws.Range("E32").Value = "'  +1.01%  "

$ws.Range("D33").Value = "'0.05706"
$ws.Range("E33").Value = "'  +3.44%  "

$ws.Range("D35").Value = "'1.304"
$ws.Range("E35").Value = "'  +3.74%  "

$ws.Range("D36").Value = "'0.7482"
$ws.Range("E36").Value = "'  +2.25%  "

$ws.Range("D37").Value = "'0.9994"
$ws.Range("E37").Value = "'  -0.04%  "

$ws.Range("E38").Value = "'  +0.10%  "

$ws.Range("D39").Value = "'0.01962"
$ws.Range("E39").Value = "'  +1.90%  "

$ws.Range("E40").Value = "'  +1.30%  "

$ws.Range("D41").Value = "'0.4510"
$ws.Range("E41").Value = "'  +2.24%  "

$ws.Range("D42").Value = "'75.02"
$ws.Range("E42").Value = "'  +3.88%  "

$ws.Range("D43").Value = "'5.993"
$ws.Range("E43").Value = "'  -1.80%  "

$ws.Range("D44").Value = "'0.8553"
$ws.Range("E44").Value = "'  +2.14%  "

$ws.Range("D45").Value = "'1.928"
$ws.Range("E45").Value = "'  +3.00%  "

$ws.Range("D46").Value = "'1.000"
$ws.Range("E46").Value = "'  -0.07%  "

$ws.Range("D47").Value = "'1.038.81"
$ws.Range("E47").Value = "'  +5.13%  "

$ws.Range("E48").Value = "'  +2.51%  "

$ws.Range("B49").Value = "'Aptos"
$ws.Range("C49").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "'7.664"
$ws.Range("E49").Value = "'  +1.79%  "

$ws.Range("B50").Value = "'EnergySwap"
$ws.Range("C50").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.946"
$ws.Range("E50").Value = "'  +2.17%  "

$ws.Range("B51").Value = "'SynthetixNetwork"
$ws.Range("C51").Value = "'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D51").Value = "'3.091"
$ws.Range("E51").Value = "'  +10.92%  "
